# Update the "Shopenzer Testcases" sheet with the new Team ID and the new
# "Executed By" names for the sprint 1 rows, then leave the workbook with
# that sheet active/selected (matching the saved view state of the
# uploaded file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shopenzer Testcases")

# New team id (PNT2022TMID12210 -> PNT2022TMID53380)
$ws.Range("F2").Value = "PNT2022TMID53380"

# New "Executed By" names for Sprint_01_TC_01 .. Sprint_01_TC_04
$ws.Range("N6").Value = "Ritunjay M"
$ws.Range("N7").Value = "Praveen Raagul R"
$ws.Range("N8").Value = "Pradeep V"
$ws.Range("N9").Value = "Munish Kumar S"

# Select the sheet / cell / zoom that was active when the file was saved.
$ws.Range("L9").Select()
$excel.ActiveWindow.Zoom = 90
